$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D1").Value = "sum"

for ($r = 2; $r -le 17; $r++) {
    $year = $ws.Cells.Item($r, 2).Text
    $split = $ws.Cells.Item($r, 3).Text
    $ws.Cells.Item($r, 4).Value = $year + " " + $split
}

$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 2023
$ws.Cells.Item(18, 3).Value = "Summer"
$ws.Cells.Item(18, 4).Value = "2023 Summer"
